$d = $word.ActiveDocument

# The paragraph currently reads (across several runs, split by a "_GoBack"
# bookmark inserted by Word at the last edit position):
#   ... files (DOCX, DOC, PDF, HTML, XPS, R | <bookmark> | TF and TXT) ...
# We need to end up with a single bold run containing the full
# "DOCX, DOC, PDF, HTML, XPS, RTF and TXT" text and no bookmark.

# Locate the full phrase (Word's Range.Text reads contiguous text even
# though it is split across runs/bookmark under the hood).
$full = $d.Content
$found = $full.Find.Execute("XPS, RTF and TXT", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Narrow the range down to just the part that straddles the bookmark
    # ("R" + bookmark + "TF and TXT"), leaving "XPS, " alone.
    $narrow = $d.Range($full.Start + 5, $full.End)

    # Replacing the text in-place (even with identical text) makes Word
    # delete the old runs/bookmark and insert a single new run, merging
    # the formatting and dropping the now-empty "_GoBack" bookmark.
    $narrow.Find.Execute("RTF and TXT", $false, $false, $false, $false, $false, $true, 1, $false, "RTF and TXT", 2) | Out-Null
}
